$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.355.69'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.231.93'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.09%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.43'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.628'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.16'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.65%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +7.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.93'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +5.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '59.52'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0947'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.16%  '

$ws.Range("E13").Value = '  -0.08%  '

$ws.Range("E14").Value = '  -0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.572.97'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.37%  '

$ws.Range("E16").Value = '  +0.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.65'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.237.68'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.25%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.314.74'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0965'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.11%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.19'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.57'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.17'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.93'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.40'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.69%  '

$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.68'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.20'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.45'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.65'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.12'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +11.97%  '

$ws.Range("E34").Value = '  +3.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0778'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.123'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.74'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.11'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0317'
$ws.Range("D40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.25'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.55'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.73'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.03'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.54'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.68%  '

$ws.Range("E46").Value = '  -1.42%  '

$ws.Range("E47").Value = '  +0.72%  '

$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("E49").Value = '  -0.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.16'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.51%  '

$ws.Range("E51").Value = '  +0.33%  '
